$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 173.4
$ws.Range("I5").Value = 173.4
$ws.Range("K5").Value = 173.4
$ws.Range("M5").Value = -58.40000000000001
$ws.Range("H74").Value = 4699.9
$ws.Range("J74").Value = 5333.3335
$ws.Range("L74").Value = 5333.3335
$ws.Range("N74").Value = -7205.3335
$ws.Range("H77").Value = 4699.9
$ws.Range("J77").Value = 5333.3335
$ws.Range("L77").Value = 26666.6675
$ws.Range("N77").Value = -36026.6675
$ws.Range("H94").Value = 20004076
$ws.Range("I94").Value = 20004076
$ws.Range("K94").Value = 20004076
$ws.Range("M94").Value = -20003625
$ws.Range("H98").Value = 1122.25
$ws.Range("I98").Value = 1122.25
$ws.Range("K98").Value = 1122.25
$ws.Range("M98").Value = 375.75
$ws.Range("H122").Value = 1122.25
$ws.Range("I122").Value = 1122.25
$ws.Range("K122").Value = 3366.75
$ws.Range("M122").Value = -916.75
$ws.Range("H132").Value = 1442.2106
$ws.Range("I132").Value = 966.80554
$ws.Range("K132").Value = 2900.41662
$ws.Range("M132").Value = -370.41662
$ws.Range("H136").Value = 79986.836
$ws.Range("J136").Value = 79986.836
$ws.Range("L136").Value = 79986.836
$ws.Range("N136").Value = -90186.836

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 39975.5
$ws.Range("J60").Value = 39975.5
$ws.Range("L60").Value = 39975.5
$ws.Range("N60").Value = -41441.5
$ws.Range("H102").Value = 2563.25
$ws.Range("I102").Value = 2563.25
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2563.25
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -941.25
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 1586.238
$ws.Range("I122").Value = 1540.55
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4621.65
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2171.65
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1253.1538
$ws.Range("I107").Value = 1365.381
$ws.Range("J107").Value = 781.8
$ws.Range("K107").Value = 1365.381
$ws.Range("L107").Value = 781.8
$ws.Range("M107").Value = 554.6189999999999
$ws.Range("N107").Value = -4621.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1874.6471
$ws.Range("I16").Value = 1797.5
$ws.Range("K16").Value = 1797.5
$ws.Range("M16").Value = -1510.5
$ws.Range("H22").Value = 205
$ws.Range("I22").Value = 205.55556
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 205.55556
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 144.44444
$ws.Range("N22").Value = -900
$ws.Range("H45").Value = 40786.75
$ws.Range("J45").Value = 40786.75
$ws.Range("L45").Value = 40786.75
$ws.Range("N45").Value = -41972.75
$ws.Range("H113").Value = 1874.6471
$ws.Range("I113").Value = 1797.5
$ws.Range("K113").Value = 1797.5
$ws.Range("M113").Value = 372.5
$ws.Range("H122").Value = 12749.5
$ws.Range("J122").Value = 13000
$ws.Range("L122").Value = 39000
$ws.Range("N122").Value = -43900
$ws.Range("H132").Value = 1795.9131
$ws.Range("I132").Value = 1785.4
$ws.Range("K132").Value = 5356.200000000001
$ws.Range("M132").Value = -2826.200000000001
$ws.Range("H134").Value = 1875.6842
$ws.Range("I134").Value = 1650.7715
$ws.Range("K134").Value = 4952.3145
$ws.Range("M134").Value = -2417.3145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2999.8572
$ws.Range("I64").Value = 1666.3334
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 4999.0002
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = -4729.0002
$ws.Range("N64").Value = -12540
$ws.Range("H67").Value = 2999.8572
$ws.Range("I67").Value = 1666.3334
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 4999.0002
$ws.Range("L67").Value = 12000
$ws.Range("M67").Value = -4063.0002
$ws.Range("N67").Value = -13872
$ws.Range("H75").Value = 9470.75
$ws.Range("J75").Value = 9916.200000000001
$ws.Range("L75").Value = 29748.6
$ws.Range("N75").Value = -31744.6
$ws.Range("H78").Value = 9470.75
$ws.Range("J78").Value = 9916.200000000001
$ws.Range("L78").Value = 89245.8
$ws.Range("N78").Value = -99229.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 9999
$ws.Range("I19").Value = 9998
$ws.Range("J19").Value = 9999.5
$ws.Range("K19").Value = 9998
$ws.Range("L19").Value = 9999.5
$ws.Range("M19").Value = -9710
$ws.Range("N19").Value = -10575.5
$ws.Range("H55").Value = 15424.5
$ws.Range("I55").Value = 7566
$ws.Range("K55").Value = 7566
$ws.Range("M55").Value = -7239
$ws.Range("H62").Value = 28056.666
$ws.Range("J62").Value = 28056.666
$ws.Range("L62").Value = 28056.666
$ws.Range("N62").Value = -29428.666
$ws.Range("H65").Value = 28056.666
$ws.Range("J65").Value = 28056.666
$ws.Range("L65").Value = 84169.99800000001
$ws.Range("N65").Value = -91033.99800000001
$ws.Range("H80").Value = 4997.8
$ws.Range("I80").Value = 3499.5
$ws.Range("J80").Value = 5996.6665
$ws.Range("K80").Value = 3499.5
$ws.Range("L80").Value = 5996.6665
$ws.Range("M80").Value = -2501.5
$ws.Range("N80").Value = -7992.6665
$ws.Range("H83").Value = 4997.8
$ws.Range("I83").Value = 3499.5
$ws.Range("J83").Value = 5996.6665
$ws.Range("K83").Value = 17497.5
$ws.Range("L83").Value = 29983.3325
$ws.Range("M83").Value = -12505.5
$ws.Range("N83").Value = -39967.3325
$ws.Range("H102").Value = 2128.5
$ws.Range("I102").Value = 2194.8823
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 2194.8823
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = -572.8823000000002
$ws.Range("N102").Value = -4244
$ws.Range("H107").Value = 35716770
$ws.Range("I107").Value = 767.6667
$ws.Range("J107").Value = 45457496
$ws.Range("K107").Value = 767.6667
$ws.Range("L107").Value = 45457496
$ws.Range("M107").Value = 1152.3333
$ws.Range("N107").Value = -45461336
$ws.Range("H113").Value = 4118.6
$ws.Range("J113").Value = 5995
$ws.Range("L113").Value = 5995
$ws.Range("N113").Value = -10335
$ws.Range("H122").Value = 2186.36
$ws.Range("I122").Value = 2236.5278
$ws.Range("K122").Value = 6709.5834
$ws.Range("M122").Value = -4259.5834
$ws.Range("H132").Value = 2246.1738
$ws.Range("I132").Value = 2078.0527
$ws.Range("J132").Value = 3044.75
$ws.Range("K132").Value = 6234.158100000001
$ws.Range("L132").Value = 9134.25
$ws.Range("M132").Value = -3704.158100000001
$ws.Range("N132").Value = -14194.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2398.2
$ws.Range("I40").Value = 2103.6316
$ws.Range("K40").Value = 2103.6316
$ws.Range("M40").Value = -1967.6316
$ws.Range("H61").Value = 18453.21
$ws.Range("I61").Value = 898.7059
$ws.Range("K61").Value = 898.7059
$ws.Range("M61").Value = -696.7059
$ws.Range("H62").Value = 20749.5
$ws.Range("J62").Value = 29500
$ws.Range("L62").Value = 29500
$ws.Range("N62").Value = -30748
$ws.Range("H65").Value = 20749.5
$ws.Range("J65").Value = 29500
$ws.Range("L65").Value = 88500
$ws.Range("N65").Value = -94740
$ws.Range("H113").Value = 18453.21
$ws.Range("I113").Value = 898.7059
$ws.Range("K113").Value = 898.7059
$ws.Range("M113").Value = 1271.2941

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 10714.286
$ws.Range("I29").Value = 10714.286
$ws.Range("K29").Value = 10714.286
$ws.Range("M29").Value = -10424.286
$ws.Range("H107").Value = 19231592
$ws.Range("I107").Value = 699.4545000000001
$ws.Range("J107").Value = 125001496
$ws.Range("K107").Value = 2098.3635
$ws.Range("L107").Value = 375004488
$ws.Range("M107").Value = -178.3635000000004
$ws.Range("N107").Value = -375008328
$ws.Range("H122").Value = 799
$ws.Range("I122").Value = 798.5
$ws.Range("K122").Value = 2395.5
$ws.Range("M122").Value = 54.5
$ws.Range("H126").Value = 2488.6785
$ws.Range("I126").Value = 2216.8696
$ws.Range("K126").Value = 6650.6088
$ws.Range("M126").Value = -4180.6088
$ws.Range("H132").Value = 15332.823
$ws.Range("I132").Value = 10373.308
$ws.Range("J132").Value = 31451.25
$ws.Range("K132").Value = 31119.924
$ws.Range("L132").Value = 94353.75
$ws.Range("M132").Value = -28589.924
$ws.Range("N132").Value = -99413.75
